# Quarterly indexing esoteric bug-fix operation
#
# Column A holds a date per row that marks the "as-of" period for the
# naive QoQ forecast on that row. These dates were incorrectly stamped on
# the 1st of the quarter-start month; the fix re-stamps them on the 15th
# of the following month (i.e. the true mid-point of the quarter) while
# leaving every other cell/value on the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = [double]$cell.Value2

    $oldDate = $epoch.AddDays($serial)
    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $newDate.ToOADate()
}
